# Update the "generate date" / "handoff" / "handback" timestamp values
# on each sheet, per the new handback report run.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file entry.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 09:02:58"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# for the first file entry.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-19 09:02:54"
$wsZhCn.Range("K2").Value = "2016-08-19 09:03:17"

# de-de sheet: "Correspond Handback DateTime" for the first file entry.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-19 09:03:23"
